# Needle calibration curvature data: re-sort the captured samples by
# timestamp (column A), ascending. The header row (1) and the very first
# sample row (row 2) stay put; rows 3-12 are the ones that get reordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A3:D12")
$sortKey   = $ws.Range("A3:A12")

$dataRange.Sort($sortKey, 1)
